# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

# ---- Rushing sheet ----
$rushing = $wb.Worksheets.Item("Rushing")

# Row 2 - D.Prescott
$rushing.Range("C2").Value = 14
$rushing.Range("D2").Value = 7
$rushing.Range("E2").Value = 16
$rushing.Range("F2").Value = 12

# Row 4 - E.Elliott
$rushing.Range("C4").Value = 135
$rushing.Range("D4").Value = 67
$rushing.Range("E4").Value = 30
$rushing.Range("F4").Value = 39

# Row 5 - T.Pollard
$rushing.Range("C5").Value = 80
$rushing.Range("D5").Value = 45
$rushing.Range("F5").Value = 17

# Row 8 - J.Hardy
$rushing.Range("D8").Value = 4

# ---- Receiving sheet ----
$receiving = $wb.Worksheets.Item("Receiving")

# Row 2 - E.Elliott
$receiving.Range("C2").Value = 51
$receiving.Range("D2").Value = 35

# Row 3 - T.Pollard
$receiving.Range("C3").Value = 39
$receiving.Range("D3").Value = 31

# Row 5 - A.Cooper
$receiving.Range("C5").Value = 82
$receiving.Range("D5").Value = 56
$receiving.Range("E5").Value = 26
$receiving.Range("F5").Value = 14
$receiving.Range("G5").Value = 21
$receiving.Range("H5").Value = 16

# Row 6 - D.Schultz
$receiving.Range("C6").Value = 84
$receiving.Range("D6").Value = 57
$receiving.Range("E6").Value = 36
$receiving.Range("G6").Value = 12
$receiving.Range("H6").Value = 6

# Row 8 - C.Lamb
$receiving.Range("C8").Value = 35
$receiving.Range("D8").Value = 27
$receiving.Range("E8").Value = 14
$receiving.Range("F8").Value = 8

# Row 10 - N.Brown
$receiving.Range("C10").Value = 12
$receiving.Range("D10").Value = 11

# Row 11 - M.Turner
$receiving.Range("C11").Value = 14

# Row 12 - B.Jarwin
$receiving.Range("C12").Value = 83
$receiving.Range("D12").Value = 67
$receiving.Range("E12").Value = 11
$receiving.Range("F12").Value = 6
